$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2363.6875
$ws.Range("I17").Value = 998
$ws.Range("J17").Value = 2678.8462
$ws.Range("K17").Value = 2994
$ws.Range("L17").Value = 8036.5386
$ws.Range("M17").Value = -2826
$ws.Range("N17").Value = -8372.5386
$ws.Range("H34").Value = 6750
$ws.Range("I34").Value = 2333.3333
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 2333.3333
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -2130.3333
$ws.Range("N34").Value = -20406
$ws.Range("H36").Value = 6750
$ws.Range("I36").Value = 2333.3333
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 2333.3333
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -1618.3333
$ws.Range("N36").Value = -21430
$ws.Range("H100").Value = 4757.8
$ws.Range("I100").Value = 3766.3333
$ws.Range("J100").Value = 6245
$ws.Range("K100").Value = 3766.3333
$ws.Range("L100").Value = 6245
$ws.Range("M100").Value = -3225.3333
$ws.Range("N100").Value = -7327
$ws.Range("H105").Value = 21747.5
$ws.Range("J105").Value = 21747.5
$ws.Range("L105").Value = 21747.5
$ws.Range("N105").Value = -28735.5
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 6385
$ws.Range("I116").Value = 12987.5
$ws.Range("K116").Value = 12987.5
$ws.Range("M116").Value = -9545.5
$ws.Range("H131").Value = 2845.2
$ws.Range("I131").Value = 2839.6667
$ws.Range("K131").Value = 8519.000100000001
$ws.Range("M131").Value = -3479.000100000001
$ws.Range("H132").Value = 4811.5
$ws.Range("I132").Value = 1307.6471
$ws.Range("J132").Value = 24666.666
$ws.Range("K132").Value = 3922.9413
$ws.Range("L132").Value = 73999.99800000001
$ws.Range("M132").Value = -1392.9413
$ws.Range("N132").Value = -79059.99800000001
$ws.Range("H137").Value = 1963.6
$ws.Range("I137").Value = 1281.5714
$ws.Range("K137").Value = 3844.7142
$ws.Range("M137").Value = -1294.7142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 13150.4
$ws.Range("I22").Value = 20005.334
$ws.Range("J22").Value = 2868
$ws.Range("K22").Value = 20005.334
$ws.Range("L22").Value = 2868
$ws.Range("M22").Value = -19706.334
$ws.Range("N22").Value = -3466

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4783.8
$ws.Range("I20").Value = 7316.6665
$ws.Range("K20").Value = 7316.6665
$ws.Range("M20").Value = -7069.6665
$ws.Range("H94").Value = 490.2
$ws.Range("I94").Value = 483.66666
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 483.66666
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = -32.66665999999998
$ws.Range("N94").Value = -1402
$ws.Range("H106").Value = 15225
$ws.Range("J106").Value = 15225
$ws.Range("L106").Value = 15225
$ws.Range("N106").Value = -17749
$ws.Range("H107").Value = 8999.625
$ws.Range("I107").Value = 1997
$ws.Range("K107").Value = 1997
$ws.Range("M107").Value = -77

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3209.875
$ws.Range("I22").Value = 1143
$ws.Range("J22").Value = 4450
$ws.Range("K22").Value = 1143
$ws.Range("L22").Value = 4450
$ws.Range("M22").Value = -793
$ws.Range("N22").Value = -5150
$ws.Range("H55").Value = 14517.75
$ws.Range("J55").Value = 19999
$ws.Range("L55").Value = 19999
$ws.Range("N55").Value = -20629
$ws.Range("H99").Value = 2945.8333
$ws.Range("J99").Value = 3807
$ws.Range("L99").Value = 3807
$ws.Range("N99").Value = -6803
$ws.Range("H107").Value = 770.0606
$ws.Range("I107").Value = 738.5263
$ws.Range("J107").Value = 812.8570999999999
$ws.Range("K107").Value = 738.5263
$ws.Range("L107").Value = 812.8570999999999
$ws.Range("M107").Value = 1181.4737
$ws.Range("N107").Value = -4652.8571
$ws.Range("H126").Value = 2945.8333
$ws.Range("J126").Value = 3807
$ws.Range("L126").Value = 11421
$ws.Range("N126").Value = -16361
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 200000
$ws.Range("J137").Value = 200000
$ws.Range("L137").Value = 200000
$ws.Range("N137").Value = -210200
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2442.5715
$ws.Range("I5").Value = 2599
$ws.Range("K5").Value = 7797
$ws.Range("M5").Value = -7685
$ws.Range("H21").Value = 54.25
$ws.Range("I21").Value = 54.25
$ws.Range("K21").Value = 162.75
$ws.Range("M21").Value = 10.25
$ws.Range("H56").Value = 6663.75
$ws.Range("I56").Value = 6663.75
$ws.Range("K56").Value = 6663.75
$ws.Range("M56").Value = -6133.75
$ws.Range("H68").Value = 983.3333
$ws.Range("J68").Value = 975
$ws.Range("L68").Value = 2925
$ws.Range("N68").Value = -4547
$ws.Range("H71").Value = 983.3333
$ws.Range("J71").Value = 975
$ws.Range("L71").Value = 8775
$ws.Range("N71").Value = -16887
$ws.Range("H113").Value = 1394.2667
$ws.Range("I113").Value = 918
$ws.Range("J113").Value = 1632.4
$ws.Range("K113").Value = 2754
$ws.Range("L113").Value = 4897.200000000001
$ws.Range("M113").Value = -584
$ws.Range("N113").Value = -9237.200000000001
$ws.Range("H135").Value = 2442.5715
$ws.Range("I135").Value = 2599
$ws.Range("K135").Value = 23391
$ws.Range("M135").Value = -20856

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872
$ws.Range("H132").Value = 51576.383
$ws.Range("I132").Value = 86154.836
$ws.Range("J132").Value = 5471.778
$ws.Range("K132").Value = 258464.508
$ws.Range("L132").Value = 16415.334
$ws.Range("M132").Value = -255934.508
$ws.Range("N132").Value = -21475.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2669.0715
$ws.Range("I7").Value = 2418.0833
$ws.Range("J7").Value = 4175
$ws.Range("K7").Value = 2418.0833
$ws.Range("L7").Value = 4175
$ws.Range("M7").Value = -2306.0833
$ws.Range("N7").Value = -4399
$ws.Range("H100").Value = 7416.6665
$ws.Range("I100").Value = 4750
$ws.Range("K100").Value = 4750
$ws.Range("M100").Value = -4209
$ws.Range("H122").Value = 3825.125
$ws.Range("I122").Value = 4079.4
$ws.Range("K122").Value = 12238.2
$ws.Range("M122").Value = -9788.200000000001
$ws.Range("H126").Value = 2669.0715
$ws.Range("I126").Value = 2418.0833
$ws.Range("J126").Value = 4175
$ws.Range("K126").Value = 7254.249899999999
$ws.Range("L126").Value = 12525
$ws.Range("M126").Value = -4784.249899999999
$ws.Range("N126").Value = -17465
$ws.Range("H132").Value = 8239
$ws.Range("I132").Value = 6126.4287
$ws.Range("K132").Value = 18379.2861
$ws.Range("M132").Value = -15849.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 65030
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H74").Value = 45000
$ws.Range("J74").Value = 45000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -46872
$ws.Range("H77").Value = 45000
$ws.Range("J77").Value = 45000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -144360
$ws.Range("H136").Value = 3225.348
$ws.Range("I136").Value = 2988.5789
$ws.Range("J136").Value = 4350
$ws.Range("K136").Value = 8965.736699999999
$ws.Range("L136").Value = 13050
$ws.Range("M136").Value = -6415.736699999999
$ws.Range("N136").Value = -18150
